$d = $word.ActiveDocument

# Locate the blank paragraph immediately preceding "Ver no Jupiter..."
# and the paragraph following "© 2020 ... Creative Commons Attribution".
# Delete the whole run of paragraphs (including their paragraph marks)
# in between, i.e. from the start of the blank paragraph up to (but not
# including) the start of the paragraph that follows the "© 2020" one.
$startRange = $null
$endRange = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $startRange = $p.Previous().Range.Start
    }
    if ($t -like "*Powered by Jekyll and Github pages*") {
        $endRange = $p.Next().Range.Start
    }
}

if ($startRange -ne $null -and $endRange -ne $null) {
    $range = $d.Range($startRange, $endRange)
    $range.Delete()
}
